$wb = $excel.ActiveWorkbook

# The localization status moved from "In Translation" to "Ready for handoff",
# and the handoff file generation timestamps were refreshed. This touches the
# "Overview" sheet (columns E/F status + G generate-date) as well as the
# per-language "zh-cn" / "de-de" sheets (column C status + column H handoff
# datetime).

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 12:58:33"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 12:58:28"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 12:58:33"

# --- Column width changes ---
# The longer "Ready for handoff" status text needs a wider status column, so
# the report generator widened the relevant columns accordingly:
#   Overview: columns E and F (status columns, 5 and 6)
#   zh-cn / de-de: column C (status column, 3)
$newStatusColWidth = 98 / 6

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
